# Sofya issues by mail on 2023-08-03
#
# Adds a new translation-key row ("notes") to the "04_ra" translations
# sheet: key column (A47) and English column (D47) both get the new
# shared string "notes" (row 47 was previously blank). Finally moves the
# active selection down to A48, matching where Excel would leave the
# cursor after typing into D47 and pressing Enter/moving to the next row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "notes"
$ws.Range("D47").Value = "notes"

$ws.Range("A48").Select() | Out-Null
